$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Tree - Wikipedia"
$ws.Range("B2").Value = "tokens@tffin`n"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 22
$ws.Range("F2").Value = "https://en.wikipedia.org/wiki/Tree"

# Row 3
$ws.Range("A3").Value = "ss-standard-user"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "https://www.arborday.org/trees/treeGuide/"

# Row 4
$ws.Range("A4").Value = "tree | Structure, Uses, Importance, & Facts | Britannica"
$ws.Range("B4").Value = "1791028@640x480`n"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "https://www.britannica.com/plant/tree"

# Row 5
$ws.Range("A5").Value = "Anniversary-logo-white"
$ws.Range("B5").Value = "info@trees.org`ninfo@trees.org`n"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "https://trees.org/"

# Row 6
$ws.Range("A6").Value = "Visa"
$ws.Range("B6").Value = "hello@onetreeplanted.org`nhello@onetreeplanted.org`nhello@onetreeplanted.org`nhello@onetreeplanted.org`n"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = "https://onetreeplanted.org/pages/why-trees"

# Row 7
$ws.Range("A7").Value = "Top 22 Benefits of Trees | TreePeople"
$ws.Range("B7").Value = "nr@context`nnr@id`n20@treepeople`n20@treepeople`n"
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "https://www.treepeople.org/tree-benefits"

# Row 8
$ws.Range("A8").Value = "https://www.youtube.com/watch%3Fv%3DHPJKxAhLw5I"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "https://www.youtube.com/watch%3Fv%3DHPJKxAhLw5I"
